# Apply "Running pipe and updated log" update to the results log worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# --- Update existing rows: the pipe run name column (B) is consolidated down
# to the two latest pipe runs (Pipe_29-08 / Pipe_29-09), and the now-finished
# NL results / Kriegstein SingleR rows are tagged with the latest run too.
$ws.Range("B83").Value = "Pipe_29-08"
$ws.Range("B84").Value = "Pipe_29-08"
$ws.Range("B89").Value = "Pipe_29-09"
$ws.Range("B90").Value = "Pipe_29-09"
$ws.Range("B91").Value = "Pipe_29-09"
$ws.Range("B92").Value = "Pipe_29-09"
$ws.Range("B93").Value = "Pipe_29-09"
$ws.Range("B94").Value = "Pipe_29-09"

# --- Append new log rows (95-102) for the KS visualization and pseudotime
# results of the latest run.
$newRows = @(
    @("results/Kriegstein", "Pipe_29-09", "KS visualization", "A+C", "oldSelection"),
    @("results/Kriegstein", "Pipe_29-09", "KS visualization", "A+C", "newSelection"),
    @("results/Kriegstein", "Pipe_29-09", "KS visualization", "N+C", "oldSelection"),
    @("results/Kriegstein", "Pipe_29-09", "KS visualization", "N+C", "newSelection"),
    @("results", "Pipe_29-09", "pseudotime", "A+C", "oldSelection"),
    @("results", "Pipe_29-09", "pseudotime", "A+C", "newSelection"),
    @("results", "Pipe_29-09", "pseudotime", "N+C", "oldSelection"),
    @("results", "Pipe_29-09", "pseudotime", "N+C", "newSelection")
)

$startRow = 95
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
}

# --- Restore the view state (active selection) to match the saved workbook.
[void]$ws.Range("E105").Select()
